$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The original sheet has data through column BD. The edit inserts three new
# columns right after column BA (i.e. before the old column BB), shifting
# everything that used to live in BB:BD three columns to the right
# (new BE:BG). The three freshly inserted columns (new BB:BD) then get their
# own values describing a new "$ghost{ string }" construct, plus a few
# unrelated one-off numeric edits lower in the sheet (BA7, BC8, BC9, BC10).
# ---------------------------------------------------------------------------

# Insert 3 new columns before the old column BB (54). This shifts the old
# BB:BD (and everything further right, e.g. BE/BF used only for col widths)
# three places to the right, exactly like Excel's "Insert Sheet Columns".
$ws.Range("BB1:BD1").EntireColumn.Insert()

# --- new column widths for the freshly inserted BB:BD -----------------------
# (Range.ColumnWidth here is quantized to a 7px "MDW" grid by the host, so we
# pick the input that lands closest to the widths produced by the edit.)
$ws.Range("BB1").ColumnWidth = 7.142857142857143
$ws.Range("BC1").ColumnWidth = 5.571428571428571
$ws.Range("BD1").ColumnWidth = 1.2857142857142858

# --- values that live in the new columns ------------------------------------
$ws.Range("BB5").Value = '$ghost{'
$ws.Range("BC5").Value = 'string'
$ws.Range("BD5").Value = '}'
$ws.Range("BC6").Value = 1
$ws.Range("BC8").Value = 2
$ws.Range("BC9").Value = 'rr'
$ws.Range("BC10").Value = 3

# --- unrelated one-off numeric edit on existing column BA -------------------
$ws.Range("BA7").Value = 1

# --- selection / view ---------------------------------------------------
$ws.Range("BC14").Select()
